$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 253, shifting rows 253:324 down to 254:325
$ws.Rows.Item(253).Insert()

# Populate the new row 253 with data (copy pattern from neighboring rows, new date/price values)
$ws.Range("A253").Value = 5
$ws.Range("B253").Value = "Macroferia Regional de Talca"
$ws.Range("C253").Value = "Maule"
$ws.Range("D253").Value = 44841
$ws.Range("E253").Value = 7
$ws.Range("F253").Value = 100112009
$ws.Range("G253").Value = "Acelga"
$ws.Range("H253").Value = "Sin especificar"
$ws.Range("I253").Value = "Primera"
$ws.Range("J253").Value = 300
$ws.Range("K253").Value = 2500
$ws.Range("L253").Value = 2500
$ws.Range("M253").Value = 2500
$ws.Range("N253").Value = "`$/docena de atados (4 kilos)"
$ws.Range("O253").Value = "Región del Maule"
$ws.Range("P253").Value = 625
$ws.Range("Q253").Value = 4
$ws.Range("R253").Value = "Hortaliza"

# Match the date-style formatting used by the rest of column D
$ws.Range("D253").NumberFormat = "YYYY-MM-DD HH:MM:SS"
